$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete rows 2-4 entirely (they are removed in the target state)
$ws.Range("A2:F4").EntireRow.Delete()

# Update row 1 with the new "international user" transaction
$ws.Range("A1").Value = "InterMasterTester"
$ws.Range("B1").Value = "USD"

# C1 / D1 hold numeric-looking text ("0" / "125") that must stay text (shared
# string), matching how the rest of the sheet stores numbers as strings.
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "0"
$ws.Range("C1").Style = "Normal"

$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Value = "125"
$ws.Range("D1").Style = "Normal"

# E1 ("American Express") is unchanged.

# F1 is removed in the new layout.
$ws.Range("F1").ClearContents()
